# Adapt column header formatting to respective input file names (#7)
#
# - Header row A1:U1 used generic "_old" / "_new" suffixes; rename them to the
#   concrete format-version suffixes "_FV2310" (old) / "_FV2404" (new).
# - Wrap the used range A1:U61 in an Excel Table ("Table1") so the renamed
#   headers double as the table's column headers.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header cells in row 1 -----------------------------------
$lastCol = 21   # A..U
for ($i = 1; $i -le $lastCol; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = $val -replace "_old$", "_FV2310"
    } elseif ($val -like "*_new") {
        $cell.Value = $val -replace "_new$", "_FV2404"
    }
}

# --- 2) Turn the used range into a Table, reusing the renamed headers ------
$lastRow = 61
$tableRange = $ws.Range("A1:U$lastRow")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
